$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value. D-column values are forced to Text format
# before assignment (and the format is reset afterwards) so that numeric-
# looking strings such as "60.94" or "1.00" keep their exact text
# representation instead of being auto-converted to numbers by Excel.
$updates = [ordered]@{
    "D2" = "60.430.88"
    "E2" = "  -0.34%  "
    "D3" = "2.626.98"
    "E3" = "  +0.40%  "
    "E4" = "  -0.05%  "
    "D5" = "522.64"
    "E5" = "  +1.27%  "
    "D6" = "151.59"
    "E6" = "  -2.00%  "
    "E7" = "  +0.13%  "
    "E8" = "  -3.95%  "
    "D9" = "6.41"
    "E9" = "  -4.44%  "
    "D11" = "0.346"
    "E11" = "  -0.31%  "
    "E12" = "  -1.00%  "
    "D13" = "3.086.85"
    "E13" = "  +0.43%  "
    "D14" = "60.408.27"
    "E14" = "  -0.39%  "
    "E15" = "  -0.67%  "
    "E16" = "  -0.40%  "
    "D17" = "2.635.47"
    "E17" = "  +0.48%  "
    "E18" = "  -1.46%  "
    "D19" = "348.22"
    "E19" = "  -2.74%  "
    "D20" = "10.48"
    "E20" = "  -1.76%  "
    "E21" = "  -0.18%  "
    "D22" = "0.995"
    "E22" = "  -0.40%  "
    "D23" = "60.94"
    "E24" = "  -0.81%  "
    "E25" = "  -0.70%  "
    "D26" = "0.992"
    "E26" = "  -0.49%  "
    "D27" = "0.0₃0839"
    "E27" = "  -0.79%  "
    "D28" = "7.18"
    "E28" = "  -2.33%  "
    "E29" = "  +0.09%  "
    "D30" = "6.09"
    "E30" = "  +2.61%  "
    "D31" = "1.61"
    "E31" = "  +1.24%  "
    "D32" = "19.11"
    "E32" = "  -1.79%  "
    "D33" = "149.85"
    "E33" = "  -1.72%  "
    "D34" = "4.02"
    "E34" = "  -0.65%  "
    "E35" = "  -1.85%  "
    "D36" = "0.888"
    "E36" = "  -0.08%  "
    "E37" = "  +5.08%  "
    "D38" = "36.56"
    "E38" = "  +0.92%  "
    "E39" = "  -2.15%  "
    "D40" = "297.67"
    "E40" = "  +2.10%  "
    "E41" = "  -1.63%  "
    "E42" = "  +0.90%  "
    "D44" = "0.997"
    "E44" = "  +0.12%  "
    "D45" = "0.0555"
    "E45" = "  -0.47%  "
    "D46" = "19.72"
    "E46" = "  +0.05%  "
    "E47" = "  +0.49%  "
    "D48" = "4.79"
    "E48" = "  -3.72%  "
    "E49" = "  +0.63%  "
    "D50" = "18.99"
    "E50" = "  -1.65%  "
    "D51" = "1.968.13"
    "E51" = "  -1.12%  "
}

foreach ($cell in $updates.Keys) {
    $range = $ws.Range($cell)
    if ($cell[0] -eq "D") {
        # Force text storage so numeric-looking values are not reinterpreted
        $range.NumberFormat = "@"
        $range.Value = $updates[$cell]
        # Drop the temporary Text number-format again so no stray style is left
        $range.Style = "Normal"
    } else {
        $range.Value = $updates[$cell]
    }
}
